$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot all data rows (2-25), columns A-R (1-18), before writing anything,
# since this edit is a permutation of whole rows and writes could otherwise clobber
# values that still need to be read as a source for another row.
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 18; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $rowVals
}

# Write back rows in their new (target) positions using the snapshot as source.
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(2, $c).Value2 = $snapshot[15][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(3, $c).Value2 = $snapshot[21][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(4, $c).Value2 = $snapshot[22][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(5, $c).Value2 = $snapshot[12][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(6, $c).Value2 = $snapshot[13][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(7, $c).Value2 = $snapshot[14][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(8, $c).Value2 = $snapshot[20][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(9, $c).Value2 = $snapshot[19][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(10, $c).Value2 = $snapshot[2][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(11, $c).Value2 = $snapshot[3][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(12, $c).Value2 = $snapshot[4][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(13, $c).Value2 = $snapshot[18][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(14, $c).Value2 = $snapshot[24][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(15, $c).Value2 = $snapshot[16][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(16, $c).Value2 = $snapshot[17][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(17, $c).Value2 = $snapshot[25][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(18, $c).Value2 = $snapshot[7][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(19, $c).Value2 = $snapshot[8][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(20, $c).Value2 = $snapshot[9][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(21, $c).Value2 = $snapshot[23][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(22, $c).Value2 = $snapshot[11][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(23, $c).Value2 = $snapshot[5][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(24, $c).Value2 = $snapshot[6][$c - 1] }
for ($c = 1; $c -le 18; $c++) { $ws.Cells.Item(25, $c).Value2 = $snapshot[10][$c - 1] }
